$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data range entirely (A1:E11) before writing the new, smaller table
$ws.Range("A1:E11").Clear()

# Update header row
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("C1").Value = "Distancia"
$ws.Range("D1").Value = "Tempo"

# Update data row
$ws.Range("A2").Value = 42
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 5344
$ws.Range("D2").Value = 0.1063048839569092
